# Insert a new "Betarraga" price record as the new row 45, shifting all
# subsequent rows (old 45..188) down by one (new 46..189).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44487
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 100114014
$ws.Range("G45").Value = "Betarraga"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 700
$ws.Range("L45").Value = 700
$ws.Range("M45").Value = 700
$ws.Range("N45").Value = "$/paquete 5 unidades"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 140
$ws.Range("Q45").Value = 5
$ws.Range("R45").Value = "Hortaliza"
